$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the 2019 data row (row 24)
$ws.Range("A24").Value = 2019
$ws.Range("B24").Value = 845
$ws.Range("C24").Value = 124250
$ws.Range("D24").Value = 153
$ws.Range("E24").Value = 32128
